# Apply refreshed cryptocurrency price / 1h-volume figures to Sheet1,
# matching the scheduled GitHub Actions data-refresh job output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.807.63'
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").Value = '2.120.65'
$ws.Range("E3").Value = '  +10.46%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.47'
$ws.Range("E5").Value = '  +4.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5250'
$ws.Range("E7").Value = '  +3.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4414'
$ws.Range("E8").Value = '  +8.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09100'
$ws.Range("E9").Value = '  +8.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.18'
$ws.Range("E10").Value = '  +11.47%  '
$ws.Range("E11").Value = '  +6.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.36'
$ws.Range("E12").Value = '  +5.48%  '
$ws.Range("D13").Value = '2.120.26'
$ws.Range("E13").Value = '  +10.64%  '
$ws.Range("E14").Value = '  +5.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.882'
$ws.Range("E15").Value = '  +8.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.29'
$ws.Range("E16").Value = '  +6.15%  '
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06637'
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.24'
$ws.Range("E20").Value = '  +3.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.418'
$ws.Range("E21").Value = '  +7.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '30.923.38'
$ws.Range("E23").Value = '  +2.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.14'
$ws.Range("E24").Value = '  +6.71%  '
$ws.Range("D25").Value = '2.366.53'
$ws.Range("E25").Value = '  +10.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.253'
$ws.Range("E26").Value = '  +2.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.00'
$ws.Range("E27").Value = '  +5.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.570'
$ws.Range("E28").Value = '  +13.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.48'
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.48'
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.185'
$ws.Range("E31").Value = '  +4.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1072'
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.270'
$ws.Range("E33").Value = '  +5.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.001'
$ws.Range("E34").Value = '  +5.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.545'
$ws.Range("E35").Value = '  +28.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02607'
$ws.Range("E36").Value = '  +6.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.604'
$ws.Range("E37").Value = '  +5.16%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06780'
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.593'
$ws.Range("E39").Value = '  +11.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.78'
$ws.Range("E40").Value = '  +11.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2278'
$ws.Range("E41").Value = '  +5.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6848'
$ws.Range("E42").Value = '  +5.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.259'
$ws.Range("E43").Value = '  +4.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.16'
$ws.Range("E44").Value = '  +4.85%  '
$ws.Range("E45").Value = '  +6.24%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.268'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.680'
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.286'
$ws.Range("E49").Value = '  +6.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.34'
$ws.Range("E50").Value = '  +5.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07088'
$ws.Range("E51").Value = '  +3.73%  '
